$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(-0.03692426742398824, 0.9235226024562185, 3.971520758787362, 1.992867471455983, 1.999337444223448, 147)
    3 = @(-0.03177683422839171, 0.9693390602487294, 3.736704271858642, 1.93305568255512, 1.9394478376698, 146)
    4 = @(-0.02353880723304746, 0.9831351250254978, 3.070892265452475, 1.752396149691181, 1.758311685393883, 145)
    5 = @(-0.02156877303751792, 1.012039586790341, 3.816035443232021, 1.953467543429381, 1.960166457363047, 144)
    6 = @(-0.01495571568351286, 1.051064849666794, 3.719246895078482, 1.928534908960292, 1.935255416224088, 143)
    7 = @(-0.008295795536640492, 1.080121458890243, 3.644641354874411, 1.909094380818929, 1.915834169966563, 142)
    8 = @(0.008388099687491768, 1.146257834064731, 3.828066810259369, 1.956544609831161, 1.963501790635553, 141)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
}
